# "add user list to project" — adds a new "users" column (E) to the
# "project hours" worksheet, listing the project-team members for each
# project row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("project hours")

# New header cell: copy the existing header formatting (bold, border,
# centered) from the neighboring "percentage" header (D1) onto E1.
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null
$ws.Range("E1").Value = "users"

# Per-project list of users, aligned with the existing rows 2-5
# (2016-08-101-01, 2016-10-105-01, 2016-10-103-02, 2016-10-103-01).
$ws.Range("E2").Value = "['HYUNG-JIN YOON', 'Kristian Lauszus', 'Ishaan Pakrasi', 'Arun Lakshmanan']"
$ws.Range("E3").Value = "['Harshal Maske']"
$ws.Range("E4").Value = "['Jonathan Hoff']"
$ws.Range("E5").Value = "['Usman Syed']"
